# Updates the cryptos list (prices/volume%) per the latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.119.34"
$ws.Range("E2").Value = "  -4.72%  "

$ws.Range("D3").Value = "3.275.53"
$ws.Range("E3").Value = "  -5.74%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.52"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -3.34%  "

$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.30"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -3.32%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  -2.79%  "

$ws.Range("D9").Value = "3.269.42"
$ws.Range("E9").Value = "  -5.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.188"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -8.05%  "

$ws.Range("E11").Value = "  -4.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.51"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -7.63%  "

$ws.Range("E13").Value = "  -6.51%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.63"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -5.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "637.14"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -0.84%  "

$ws.Range("D16").Value = "3.803.95"
$ws.Range("E16").Value = "  -5.76%  "

$ws.Range("D17").Value = "66.056.14"
$ws.Range("E17").Value = "  -4.47%  "

$ws.Range("E18").Value = "  -1.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.116"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -3.50%  "

$ws.Range("D20").Value = "3.277.88"
$ws.Range("E20").Value = "  -5.36%  "

$ws.Range("E21").Value = "  -7.84%  "

$ws.Range("E22").Value = "  -3.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.43"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +3.76%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "107.83"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +8.91%  "

$ws.Range("E25").Value = "  -7.21%  "

$ws.Range("E26").Value = "  -7.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.67"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -6.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.60"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -3.38%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.72"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -6.45%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.35"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -6.40%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.09"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -4.78%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.27"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -6.80%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.05"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -4.88%  "

$ws.Range("E34").Value = "  -3.75%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "536.11"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +2.26%  "

$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "3.741.38"
$ws.Range("E36").Value = "  +0.58%  "

$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "57.53"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -5.50%  "

$ws.Range("E38").Value = "  -0.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.36"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -3.62%  "

$ws.Range("D40").Value = "0.0₃0731"
$ws.Range("E40").Value = "  -7.76%  "

$ws.Range("E41").Value = "  -1.92%  "

$ws.Range("E42").Value = "  -6.84%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.46"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -1.51%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "32.86"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -4.10%  "

$ws.Range("E45").Value = "  -8.93%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.28"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -2.20%  "

$ws.Range("E47").Value = "  -6.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.62"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -7.07%  "

$ws.Range("E49").Value = "  -3.68%  "

$ws.Range("E50").Value = "  -0.05%  "

$ws.Range("E51").Value = "  +2.07%  "
